$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.086.34"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "1.821.80"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.47"
$ws.Range("E5").Value = "  -1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4627"
$ws.Range("E7").Value = "  -2.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3637"
$ws.Range("E8").Value = "  -1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8690"
$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.12"
$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("D12").Value = "1.876.25"
$ws.Range("E12").Value = "  +2.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07572"
$ws.Range("E13").Value = "  +2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.341"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.54"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.473"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008643"
$ws.Range("E18").Value = "  -2.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").Value = "27.378.69"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  -2.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.197"
$ws.Range("E22").Value = "  -2.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("D24").Value = "2.093.04"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.70"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.874"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.18"
$ws.Range("E27").Value = "  -2.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("E28").Value = "  -3.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.18"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.066"
$ws.Range("E30").Value = "  -4.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08909"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.962"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7343"
$ws.Range("E33").Value = "  -3.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.451"
$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.502"
$ws.Range("E37").Value = "  +4.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05254"
$ws.Range("E38").Value = "  -2.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01919"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.939"
$ws.Range("E41").Value = "  -2.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.151"
$ws.Range("E42").Value = "  -2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5205"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1631"
$ws.Range("E44").Value = "  -2.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.260"
$ws.Range("E45").Value = "  -3.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4883"
$ws.Range("E46").Value = "  -1.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("E48").Value = "  -3.39%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.81"
$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.633"
$ws.Range("E50").Value = "  -3.02%  "

$ws.Range("E51").Value = "  -1.33%  "
